# Add a new "Test Name" payroll record as a new sheet, and drop two
# now-obsolete payroll rows from Sheet1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Remove the last two payroll rows (Test Employee 4 & 5) from Sheet1.
$ws1.Rows.Item(5).Resize(2).Delete()

# New worksheet, placed after Sheet2, to store the new name/value pair.
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Test Name"

$ws3.Range("A1").Value = "Test Name"
$ws3.Range("A1").Font.Bold = $true

# Type a throwaway placeholder first, then overwrite it with the real
# random string (mirrors how it was actually entered).
$ws3.Range("A2").Value = "aaaaaa"
$ws3.Range("A2").Value = "uAJINIPj"

$ws3.Columns.Item(1).ColumnWidth = 22.75

$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

# Restore/refresh selections on the other sheets.
$ws1.Range("G11").Select() | Out-Null
$ws2.Range("A2").Select() | Out-Null
$ws3.Range("A2").Select() | Out-Null

$ws3.Activate()
